# Update "想去人数" (F column) values across the four sheets of the
# 广州-漫展信息 workbook, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        $SheetName,
        $Updates
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $Updates.Keys) {
        $ws.Range("F$row").Value = $Updates[$row]
    }
}

# Sheet "展览" (sheet1)
Set-FValues "展览" @{
    2  = 1434
    3  = 1410
    6  = 668
    7  = 27
    8  = 615
    10 = 74
    11 = 1364
    12 = 31814
    13 = 6861
    15 = 342
    16 = 564
    19 = 86
    21 = 428
    22 = 91
    23 = 776
    24 = 314
    26 = 423
    28 = 183
    30 = 725
    33 = 713
    34 = 105
    36 = 777
    37 = 281
    38 = 49
}

# Sheet "演出" (sheet2)
Set-FValues "演出" @{
    2  = 1144
    7  = 4315
    9  = 229
    12 = 56
    14 = 5
    19 = 4282
}

# Sheet "本地生活" (sheet3)
Set-FValues "本地生活" @{
    2 = 1421
    3 = 346
}

# Sheet "全部类型" (sheet4)
Set-FValues "全部类型" @{
    2  = 1421
    3  = 346
    4  = 1144
    5  = 1434
    6  = 1410
    8  = 668
    9  = 27
    10 = 615
    12 = 74
    13 = 1364
    16 = 229
    17 = 229
    20 = 6861
    22 = 342
    23 = 56
    24 = 564
    27 = 86
    28 = 5
    31 = 428
    32 = 91
    33 = 776
    34 = 314
    36 = 423
    38 = 183
    40 = 725
    44 = 105
    45 = 777
    46 = 281
    47 = 49
}
